# Refresh realtime "offerte" report with new data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report date shown in the I1 header cell.
$ws.Range("I1").Value = "24/03/2023"

# Updated data rows (r=2..20), columns B..J:
#   B: Attesa, C: Risposte Eff., D: Offerte, E: Abb sup. 14, F: Abb inf. 14,
#   G: Short Call min 10, H: Cleared, I: (metric), J: Delta_Offerto
$data = @(
    @(2, 179, 182, 1, 0, 3, 0, 270.4, -32.69230769230769),
    @(0, 108, 108, 0, 0, 2, 0, 159, -32.07547169811321),
    @(0, 1, 1, 0, 0, 0, 0, 2, -50),
    @(5, 83, 95, 6, 1, 0, 0, 73, 30.13698630136987),
    @(0, 13, 13, 0, 0, 0, 0, 29, -55.17241379310344),
    @(1, 33, 36, 2, 0, 3, 0, 50, -28),
    @(2, 111, 116, 3, 0, 1, 0, 29, 300),
    @(0, 14, 14, 0, 1, 0, 0, 71, -80.28169014084507),
    @(25, 129, 195, 39, 3, 0, 0, 293, -33.44709897610921),
    @(0, 68, 68, 0, 0, 4, 0, 204, -66.66666666666667),
    @(12, 134, 190, 11, 0, 3, 32, 331.5, -42.68476621417798),
    @(0, 23, 23, 0, 0, 0, 0, 233, -90.12875536480686),
    @(5, 213, 243, 5, 2, 1, 18, 354, -31.35593220338983),
    @(0, 76, 80, 4, 0, 0, 0, 113, -29.20353982300885),
    @(0, 97, 124, 25, 2, 2, 0, 151, -17.88079470198676),
    @(0, 18, 18, 0, 0, 1, 0, 47, -61.70212765957447),
    @(0, 0, 0, 0, 0, 0, 0, 1, -100),
    @(0, 3, 3, 0, 0, 0, 0, 6, -50),
    @(0, 10, 10, 0, 0, 0, 0, 28, -64.28571428571428)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        # Columns B(2) through J(10)
        $ws.Cells.Item($rowIndex, $c + 2).Value = $rowValues[$c]
    }
}
